$d = $word.ActiveDocument
$full = $d.Content
$xml = $full.WordOpenXML

# Remove every <w:contextualSpacing .../> element (self-closing, any attrs)
# wherever it occurs across the package (document.xml, comments.xml, ...).
$pattern = '<w:contextualSpacing[^>]*/>'
$newXml = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern, '')

$full.WordOpenXML = $newXml

Write-Output "done"
